$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / column B holds the "CasesTab" Cypher query. The commit removes the
# trailing "Cohort" coalesce line (and the now-dangling comma after the
# "Response to Treatment" line) from that query text.
$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Urethra, Prostate'] and diag.best_response in ['Not Determined']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesTabQuery

# The text shrank (one fewer wrapped line), so Excel recalculated the
# auto-fit row height for row 2 (rows 3/4 are untouched).
$ws.Rows.Item(2).RowHeight = 319

# Update the saved view/selection state to match the new focus cell.
$ws.Range("B2").Select()
